# Generate Report for handoff
#
# Updates the localization-status workbook so that the zh-cn and de-de
# source files show they are ready for handoff: the "Handoff transform
# failed" status becomes "Ready for handoff" everywhere it is shown
# (Overview sheet + the per-locale sheets), and each per-locale sheet
# gets its first row's handoff file / handoff datetime / handoff reason
# filled in, with the handoff file name rendered as a hyperlink.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Handoff transform failed" -> "Ready for handoff" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B2").Value = "Ready for handoff"

$zhcnFile = "c8b52b61-6c53-4fcd-b4bf-62d15daec756.53d5fc92a0a36555942465a9686cf153479819fc.zh-cn.xlf"
$zhcnUrl = "https://github.com/OpenLocalizationTest/oltest/blob/d14d00e03cdb26f629d4bb04b19e995b7d9ef649/e2e/" + $zhcnFile
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), $zhcnUrl, "", "", $zhcnFile)
$zhcn.Range("C2").Font.Underline = $true
$zhcn.Range("C2").Font.Color = 15570276

$zhcn.Range("D2").Value = "2016-02-17 03:19:48"
$zhcn.Range("H2").Value = "Include"

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B2").Value = "Ready for handoff"

$dedeFile = "c8b52b61-6c53-4fcd-b4bf-62d15daec756.53d5fc92a0a36555942465a9686cf153479819fc.de-de.xlf"
$dedeUrl = "https://github.com/OpenLocalizationTest/oltest/blob/d14d00e03cdb26f629d4bb04b19e995b7d9ef649/e2e/" + $dedeFile
$dede.Hyperlinks.Add($dede.Range("C2"), $dedeUrl, "", "", $dedeFile)
$dede.Range("C2").Font.Underline = $true
$dede.Range("C2").Font.Color = 15570276

$dede.Range("D2").Value = "2016-02-17 03:19:58"
$dede.Range("H2").Value = "Include"
